# Apply updated loading_percent values (case with 380 kV done)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 23.74456521784683
$ws.Cells.Item(2, 4).Value2 = 9.581584845702061
$ws.Cells.Item(2, 5).Value2 = 10.56676181936496
$ws.Cells.Item(2, 6).Value2 = 67.63404625294729
$ws.Cells.Item(2, 7).Value2 = 3.828652643670688
$ws.Cells.Item(2, 10).Value2 = 8.751089675392741
$ws.Cells.Item(2, 12).Value2 = 15.8261477676847
$ws.Cells.Item(2, 13).Value2 = 22.67764792265407
$ws.Cells.Item(3, 2).Value2 = 23.62782581568237
$ws.Cells.Item(3, 4).Value2 = 9.251384897604128
$ws.Cells.Item(3, 5).Value2 = 9.94356199153588
$ws.Cells.Item(3, 6).Value2 = 67.13156621350963
$ws.Cells.Item(3, 7).Value2 = 3.836024405256647
$ws.Cells.Item(3, 10).Value2 = 8.599196511442869
$ws.Cells.Item(3, 12).Value2 = 15.76632709752477
$ws.Cells.Item(3, 13).Value2 = 22.63978301522346
$ws.Cells.Item(4, 2).Value2 = 23.56302650621591
$ws.Cells.Item(4, 4).Value2 = 9.046475311849449
$ws.Cells.Item(4, 5).Value2 = 9.539646327897492
$ws.Cells.Item(4, 6).Value2 = 66.83880945373981
$ws.Cells.Item(4, 7).Value2 = 3.840771993281896
$ws.Cells.Item(4, 10).Value2 = 8.503077729561562
$ws.Cells.Item(4, 12).Value2 = 15.73463869800957
$ws.Cells.Item(4, 13).Value2 = 22.62283977775602
$ws.Cells.Item(5, 2).Value2 = 23.53836948625679
$ws.Cells.Item(5, 4).Value2 = 8.962568374976183
$ws.Cells.Item(5, 5).Value2 = 9.369732301358361
$ws.Cells.Item(5, 6).Value2 = 66.72352425234671
$ws.Cells.Item(5, 7).Value2 = 3.842762629703776
$ws.Cells.Item(5, 10).Value2 = 8.463206854985817
$ws.Cells.Item(5, 12).Value2 = 15.7230006184811
$ws.Cells.Item(5, 13).Value2 = 22.61752387789361
$ws.Cells.Item(6, 2).Value2 = 23.53438137149175
$ws.Cells.Item(6, 4).Value2 = 8.948615535782338
$ws.Cells.Item(6, 5).Value2 = 9.341198373760085
$ws.Cells.Item(6, 6).Value2 = 66.70462508044847
$ws.Cells.Item(6, 7).Value2 = 3.843096561982585
$ws.Cells.Item(6, 10).Value2 = 8.456544442568052
$ws.Cells.Item(6, 12).Value2 = 15.72114532665927
$ws.Cells.Item(6, 13).Value2 = 22.61673718307372
$ws.Cells.Item(7, 2).Value2 = 23.5626868667305
$ws.Cells.Item(7, 4).Value2 = 9.0453451596668
$ws.Cells.Item(7, 5).Value2 = 9.537376276938026
$ws.Cells.Item(7, 6).Value2 = 66.83723835747649
$ws.Cells.Item(7, 7).Value2 = 3.840798612727804
$ws.Cells.Item(7, 10).Value2 = 8.502542835835262
$ws.Cells.Item(7, 12).Value2 = 15.73447657098496
$ws.Cells.Item(7, 13).Value2 = 22.62276165057973
$ws.Cells.Item(8, 2).Value2 = 23.70289812778019
$ws.Cells.Item(8, 4).Value2 = 9.468262156683602
$ws.Cells.Item(8, 5).Value2 = 10.35631325013862
$ws.Cells.Item(8, 6).Value2 = 67.4575511998592
$ws.Cells.Item(8, 7).Value2 = 3.831148683796122
$ws.Cells.Item(8, 10).Value2 = 8.699319192579019
$ws.Cells.Item(8, 12).Value2 = 15.80448014382632
$ws.Cells.Item(8, 13).Value2 = 22.66328545393848
$ws.Cells.Item(9, 2).Value2 = 24.03146355904677
$ws.Cells.Item(9, 4).Value2 = 10.27480272887891
$ws.Cells.Item(9, 5).Value2 = 11.79256163426518
$ws.Cells.Item(9, 6).Value2 = 68.79652098787132
$ws.Cells.Item(9, 7).Value2 = 3.813966951352783
$ws.Cells.Item(9, 10).Value2 = 9.06186740468379
$ws.Cells.Item(9, 12).Value2 = 15.98135868313061
$ws.Cells.Item(9, 13).Value2 = 22.79260128214048
$ws.Cells.Item(10, 2).Value2 = 24.30408155643881
$ws.Cells.Item(10, 4).Value2 = 10.84672752465685
$ws.Cells.Item(10, 5).Value2 = 12.74367458142317
$ws.Cells.Item(10, 6).Value2 = 69.85105592517245
$ws.Cells.Item(10, 7).Value2 = 3.802385472961732
$ws.Cells.Item(10, 10).Value2 = 9.313185572820913
$ws.Cells.Item(10, 12).Value2 = 16.13480921392114
$ws.Cells.Item(10, 13).Value2 = 22.9176661488088
$ws.Cells.Item(11, 2).Value2 = 24.43452476981533
$ws.Cells.Item(11, 4).Value2 = 11.10122640144792
$ws.Cells.Item(11, 5).Value2 = 13.1537454187919
$ws.Cells.Item(11, 6).Value2 = 70.34519936348151
$ws.Cells.Item(11, 7).Value2 = 3.797338694761532
$ws.Cells.Item(11, 10).Value2 = 9.42409482693669
$ws.Cells.Item(11, 12).Value2 = 16.20955112398211
$ws.Cells.Item(11, 13).Value2 = 22.98098899419965
$ws.Cells.Item(12, 2).Value2 = 24.48481151996785
$ws.Cells.Item(12, 4).Value2 = 11.19669705337577
$ws.Cells.Item(12, 5).Value2 = 13.30578341386607
$ws.Cells.Item(12, 6).Value2 = 70.53430368967136
$ws.Cells.Item(12, 7).Value2 = 3.795459148676798
$ws.Cells.Item(12, 10).Value2 = 9.465590086313378
$ws.Cells.Item(12, 12).Value2 = 16.23854734957902
$ws.Cells.Item(12, 13).Value2 = 23.00588179518303
$ws.Cells.Item(13, 2).Value2 = 24.47394230493749
$ws.Cells.Item(13, 4).Value2 = 11.17617715569272
$ws.Cells.Item(13, 5).Value2 = 13.27318364552404
$ws.Cells.Item(13, 6).Value2 = 70.493490010668
$ws.Cells.Item(13, 7).Value2 = 3.795862544073415
$ws.Cells.Item(13, 10).Value2 = 9.456675926820743
$ws.Cells.Item(13, 12).Value2 = 16.23227194875431
$ws.Cells.Item(13, 13).Value2 = 23.0004802290501
$ws.Cells.Item(14, 2).Value2 = 24.43864423756425
$ws.Cells.Item(14, 4).Value2 = 11.10909937891622
$ws.Cells.Item(14, 5).Value2 = 13.16631872374413
$ws.Cells.Item(14, 6).Value2 = 70.36071789242982
$ws.Cells.Item(14, 7).Value2 = 3.797183432722681
$ws.Cells.Item(14, 10).Value2 = 9.427518828415176
$ws.Cells.Item(14, 12).Value2 = 16.21192287098803
$ws.Cells.Item(14, 13).Value2 = 22.98301869235826
$ws.Cells.Item(15, 2).Value2 = 24.41713812728944
$ws.Cells.Item(15, 4).Value2 = 11.06789236281832
$ws.Cells.Item(15, 5).Value2 = 13.10043815588341
$ws.Cells.Item(15, 6).Value2 = 70.27964656382618
$ws.Cells.Item(15, 7).Value2 = 3.797996615565408
$ws.Cells.Item(15, 10).Value2 = 9.409593321371108
$ws.Cells.Item(15, 12).Value2 = 16.19954821424743
$ws.Cells.Item(15, 13).Value2 = 22.97244165902861
$ws.Cells.Item(16, 2).Value2 = 24.29568329859642
$ws.Cells.Item(16, 4).Value2 = 10.82997361907974
$ws.Cells.Item(16, 5).Value2 = 12.71642091986
$ws.Cells.Item(16, 6).Value2 = 69.81904509010639
$ws.Cells.Item(16, 7).Value2 = 3.802719725540603
$ws.Cells.Item(16, 10).Value2 = 9.30586756163995
$ws.Cells.Item(16, 12).Value2 = 16.13002258428392
$ws.Cells.Item(16, 13).Value2 = 22.91365634643884
$ws.Cells.Item(17, 2).Value2 = 24.22279767914853
$ws.Cells.Item(17, 4).Value2 = 10.68250021490456
$ws.Cells.Item(17, 5).Value2 = 12.47505198629893
$ws.Cells.Item(17, 6).Value2 = 69.54011284042801
$ws.Cells.Item(17, 7).Value2 = 3.805673754835386
$ws.Cells.Item(17, 10).Value2 = 9.24135097722524
$ws.Cells.Item(17, 12).Value2 = 16.08862369015085
$ws.Cells.Item(17, 13).Value2 = 22.87923308438152
$ws.Cells.Item(18, 2).Value2 = 24.181482949625
$ws.Cells.Item(18, 4).Value2 = 10.59714800042299
$ws.Cells.Item(18, 5).Value2 = 12.33409579253035
$ws.Cells.Item(18, 6).Value2 = 69.38104245457396
$ws.Cells.Item(18, 7).Value2 = 3.807393720581702
$ws.Cells.Item(18, 10).Value2 = 9.203921830092732
$ws.Cells.Item(18, 12).Value2 = 16.06527787687703
$ws.Cells.Item(18, 13).Value2 = 22.86003970846786
$ws.Cells.Item(19, 2).Value2 = 24.16759976649367
$ws.Cells.Item(19, 4).Value2 = 10.56816109554217
$ws.Cells.Item(19, 5).Value2 = 12.28600509847016
$ws.Cells.Item(19, 6).Value2 = 69.32742096686056
$ws.Cells.Item(19, 7).Value2 = 3.807979668363181
$ws.Cells.Item(19, 10).Value2 = 9.191194239563115
$ws.Cells.Item(19, 12).Value2 = 16.05745386616771
$ws.Cells.Item(19, 13).Value2 = 22.85364554700832
$ws.Cells.Item(20, 2).Value2 = 24.2304938780053
$ws.Cells.Item(20, 4).Value2 = 10.69825446445911
$ws.Cells.Item(20, 5).Value2 = 12.50096627053782
$ws.Cells.Item(20, 6).Value2 = 69.56966513203925
$ws.Cells.Item(20, 7).Value2 = 3.805357134040585
$ws.Cells.Item(20, 10).Value2 = 9.248252183817373
$ws.Cells.Item(20, 12).Value2 = 16.09298258873571
$ws.Cells.Item(20, 13).Value2 = 22.88283485791623
$ws.Cells.Item(21, 2).Value2 = 24.44898823700969
$ws.Cells.Item(21, 4).Value2 = 11.12882688431688
$ws.Cells.Item(21, 5).Value2 = 13.19779562676601
$ws.Cells.Item(21, 6).Value2 = 70.39966322138753
$ws.Cells.Item(21, 7).Value2 = 3.796794601760962
$ws.Cells.Item(21, 10).Value2 = 9.436096738025125
$ws.Cells.Item(21, 12).Value2 = 16.21788121810454
$ws.Cells.Item(21, 13).Value2 = 22.98812286095886
$ws.Cells.Item(22, 2).Value2 = 24.59696214376301
$ws.Cells.Item(22, 4).Value2 = 11.40493690994145
$ws.Cells.Item(22, 5).Value2 = 13.63429175801279
$ws.Cells.Item(22, 6).Value2 = 70.95362737792664
$ws.Cells.Item(22, 7).Value2 = 3.791382293182824
$ws.Cells.Item(22, 10).Value2 = 9.555923212776296
$ws.Cells.Item(22, 12).Value2 = 16.30354089401507
$ws.Cells.Item(22, 13).Value2 = 23.06225515648643
$ws.Cells.Item(23, 2).Value2 = 24.51752365912644
$ws.Cells.Item(23, 4).Value2 = 11.25808226448171
$ws.Cells.Item(23, 5).Value2 = 13.40305528776704
$ws.Cells.Item(23, 6).Value2 = 70.65694440663435
$ws.Cells.Item(23, 7).Value2 = 3.794254232751707
$ws.Cells.Item(23, 10).Value2 = 9.492242479714779
$ws.Cells.Item(23, 12).Value2 = 16.25745971650353
$ws.Cells.Item(23, 13).Value2 = 23.02220643416072
$ws.Cells.Item(24, 2).Value2 = 24.22701259191325
$ws.Cells.Item(24, 4).Value2 = 10.69113373007436
$ws.Cells.Item(24, 5).Value2 = 12.48925724405764
$ws.Cells.Item(24, 6).Value2 = 69.55630051139565
$ws.Cells.Item(24, 7).Value2 = 3.805500210827351
$ws.Cells.Item(24, 10).Value2 = 9.245133199238619
$ws.Cells.Item(24, 12).Value2 = 16.09101051232221
$ws.Cells.Item(24, 13).Value2 = 22.88120463572297
$ws.Cells.Item(25, 2).Value2 = 23.93698029771536
$ws.Cells.Item(25, 4).Value2 = 10.05979205830288
$ws.Cells.Item(25, 5).Value2 = 11.42227782905385
$ws.Cells.Item(25, 6).Value2 = 68.42155246707961
$ws.Cells.Item(25, 7).Value2 = 3.818430646247676
$ws.Cells.Item(25, 10).Value2 = 8.966374792038232
$ws.Cells.Item(25, 12).Value2 = 15.92932901113056
$ws.Cells.Item(25, 13).Value2 = 22.75231014120296
